# Labels30.xlsx — "Added more healthy images and extraced features"
#
# The source workbook is a single-column label list (A1 header "Label",
# A2:A301 = 0/1 labels). This commit:
#   1) touches the existing label cells A1:A301 (re-applying their style,
#      which is what produced the second cellXfs entry — same font/number
#      format, just with the "applyFont" flag flipped on), and
#   2) appends 30 new rows (A302:A331), all valued 0, extending the
#      dataset with newly-extracted (negative/"unhealthy") labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Re-apply the Normal style across the existing data range. This is a
#    no-visual-effect restyle, matching the new cellXfs entry added in the
#    diff (identical numFmt/font, only the "applied" flags differ from the
#    original style).
$ws.Range("A1:A301").Style = "Normal"

# 2) Append the 30 new label rows at the bottom, all 0.
for ($r = 302; $r -le 331; $r++) {
    $ws.Cells.Item($r, 1).Value2 = 0
}

# 3) Leave the selection on the new last cell, like a user who just
#    finished typing/pasting the new rows.
$ws.Range("A332").Select()
